# Insert a new weekly record at row 6 (just below the existing row 5),
# which pushes the previous rows 6..48 down to 7..49 and extends the
# sheet's used range from A1:R48 to A1:R49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new record. Columns A, B, C,
# E, F, G, H, I, N, O, Q, R carry the same constant values as every other
# data row in this sheet; D, J, K, L, M, P hold the new observation.
$ws.Cells.Item(6, 1).Value = 9
$ws.Cells.Item(6, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value = "Metropolitana"
$ws.Cells.Item(6, 4).Value = 44761
$ws.Cells.Item(6, 5).Value = 13
$ws.Cells.Item(6, 6).Value = 100112035
$ws.Cells.Item(6, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 43
$ws.Cells.Item(6, 11).Value = 19000
$ws.Cells.Item(6, 12).Value = 19000
$ws.Cells.Item(6, 13).Value = 19000
$ws.Cells.Item(6, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(6, 15).Value = "Hijuelas"
$ws.Cells.Item(6, 16).Value = 1267
$ws.Cells.Item(6, 17).Value = 15
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Keep the date number format (numFmtId 165) that every other "Fecha"
# cell in column D uses.
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
